$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 300600.62
$ws.Range("I33").Value = 426.52942
$ws.Range("J33").Value = 867596.1
$ws.Range("K33").Value = 426.52942
$ws.Range("L33").Value = 867596.1
$ws.Range("M33").Value = -197.52942
$ws.Range("N33").Value = -868054.1
$ws.Range("H92").Value = 571.8261
$ws.Range("I92").Value = 517.85
$ws.Range("K92").Value = 517.85
$ws.Range("M92").Value = 730.15
$ws.Range("H121").Value = 930.2222
$ws.Range("J121").Value = 849.6
$ws.Range("L121").Value = 2548.8
$ws.Range("N121").Value = -6042.8
$ws.Range("H137").Value = 1129.258
$ws.Range("I137").Value = 1102.16
$ws.Range("K137").Value = 3306.48
$ws.Range("M137").Value = -756.4800000000005
$ws.Range("H138").Value = 2253.6667
$ws.Range("I138").Value = 1866.2667
$ws.Range("J138").Value = 2328.1667
$ws.Range("K138").Value = 5598.800099999999
$ws.Range("L138").Value = 6984.500100000001
$ws.Range("M138").Value = -458.8000999999995
$ws.Range("N138").Value = -17264.5001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5104.49
$ws.Range("I32").Value = 4695.3506
$ws.Range("J32").Value = 18333.334
$ws.Range("K32").Value = 4695.3506
$ws.Range("L32").Value = 18333.334
$ws.Range("M32").Value = -4408.3506
$ws.Range("N32").Value = -18907.334
$ws.Range("H61").Value = 2145.4285
$ws.Range("I61").Value = 1238.5333
$ws.Range("J61").Value = 3191.8462
$ws.Range("K61").Value = 1238.5333
$ws.Range("L61").Value = 3191.8462
$ws.Range("M61").Value = -1026.5333
$ws.Range("N61").Value = -3615.8462
$ws.Range("H122").Value = 2288.5
$ws.Range("I122").Value = 2352.4443
$ws.Range("J122").Value = 2096.6667
$ws.Range("K122").Value = 7057.3329
$ws.Range("L122").Value = 6290.000100000001
$ws.Range("M122").Value = -4607.3329
$ws.Range("N122").Value = -11190.0001
$ws.Range("H136").Value = 2145.4285
$ws.Range("I136").Value = 1238.5333
$ws.Range("J136").Value = 3191.8462
$ws.Range("K136").Value = 3715.5999
$ws.Range("L136").Value = 9575.5386
$ws.Range("M136").Value = -1165.5999
$ws.Range("N136").Value = -14675.5386
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 40291.406
$ws.Range("I31").Value = 586.35297
$ws.Range("K31").Value = 586.35297
$ws.Range("M31").Value = -291.35297
$ws.Range("H34").Value = 40291.406
$ws.Range("I34").Value = 586.35297
$ws.Range("K34").Value = 586.35297
$ws.Range("M34").Value = -384.35297
$ws.Range("H99").Value = 9882.166999999999
$ws.Range("I99").Value = 1440.1666
$ws.Range("J99").Value = 18324.166
$ws.Range("K99").Value = 1440.1666
$ws.Range("L99").Value = 18324.166
$ws.Range("M99").Value = 57.83339999999998
$ws.Range("N99").Value = -21320.166
$ws.Range("H126").Value = 9882.166999999999
$ws.Range("I126").Value = 1440.1666
$ws.Range("J126").Value = 18324.166
$ws.Range("K126").Value = 4320.4998
$ws.Range("L126").Value = 54972.49800000001
$ws.Range("M126").Value = -1850.4998
$ws.Range("N126").Value = -59912.49800000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 612.4167
$ws.Range("I122").Value = 619.4
$ws.Range("J122").Value = 610.5789
$ws.Range("K122").Value = 5574.599999999999
$ws.Range("L122").Value = 5495.2101
$ws.Range("M122").Value = -3124.599999999999
$ws.Range("N122").Value = -10395.2101
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 13500.875
$ws.Range("J21").Value = 13500.875
$ws.Range("L21").Value = 13500.875
$ws.Range("N21").Value = -13846.875
$ws.Range("H30").Value = 13500.875
$ws.Range("J30").Value = 13500.875
$ws.Range("L30").Value = 13500.875
$ws.Range("N30").Value = -13710.875
$ws.Range("H122").Value = 10000
$ws.Range("I122").Value = 10000
$ws.Range("K122").Value = 30000
$ws.Range("M122").Value = -27550
$ws.Range("H126").Value = 6539134.5
$ws.Range("I126").Value = 3852.4
$ws.Range("J126").Value = 14708237
$ws.Range("K126").Value = 11557.2
$ws.Range("L126").Value = 44124711
$ws.Range("M126").Value = -9087.200000000001
$ws.Range("N126").Value = -44129651
$ws.Range("H132").Value = 2832.7568
$ws.Range("I132").Value = 2250.5356
$ws.Range("K132").Value = 6751.6068
$ws.Range("M132").Value = -4221.6068
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4900
$ws.Range("I7").Value = 2500
$ws.Range("J7").Value = 6100
$ws.Range("K7").Value = 2500
$ws.Range("L7").Value = 6100
$ws.Range("M7").Value = -2388
$ws.Range("N7").Value = -6324
$ws.Range("H40").Value = 65035.75
$ws.Range("I40").Value = 251700.75
$ws.Range("J40").Value = 2814.0833
$ws.Range("K40").Value = 251700.75
$ws.Range("L40").Value = 2814.0833
$ws.Range("M40").Value = -251564.75
$ws.Range("N40").Value = -3086.0833
$ws.Range("H55").Value = 227859.27
$ws.Range("J55").Value = 463.73077
$ws.Range("L55").Value = 463.73077
$ws.Range("N55").Value = -809.73077
$ws.Range("H122").Value = 2930.4443
$ws.Range("I122").Value = 2930.4443
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8791.332900000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6341.332900000001
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 4900
$ws.Range("I126").Value = 2500
$ws.Range("J126").Value = 6100
$ws.Range("K126").Value = 7500
$ws.Range("L126").Value = 18300
$ws.Range("M126").Value = -5030
$ws.Range("N126").Value = -23240
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2672.5334
$ws.Range("I122").Value = 1852.5385
$ws.Range("J122").Value = 8002.5
$ws.Range("K122").Value = 5557.6155
$ws.Range("L122").Value = 24007.5
$ws.Range("M122").Value = -3107.6155
$ws.Range("N122").Value = -28907.5
$ws.Range("H126").Value = 1457.05
$ws.Range("I126").Value = 1814.8334
$ws.Range("J126").Value = 920.375
$ws.Range("K126").Value = 5444.5002
$ws.Range("L126").Value = 2761.125
$ws.Range("M126").Value = -2974.5002
$ws.Range("N126").Value = -7701.125
$ws.Range("H132").Value = 3155.8647
$ws.Range("I132").Value = 3467
$ws.Range("J132").Value = 2315.8
$ws.Range("K132").Value = 10401
$ws.Range("L132").Value = 6947.400000000001
$ws.Range("M132").Value = -7871
$ws.Range("N132").Value = -12007.4
$ws.Range("H136").Value = 1469.7222
$ws.Range("I136").Value = 475.23334
$ws.Range("J136").Value = 2712.8333
$ws.Range("K136").Value = 1425.70002
$ws.Range("L136").Value = 8138.499899999999
$ws.Range("M136").Value = 1124.29998
$ws.Range("N136").Value = -13238.4999
